# Harmonize parameters w Leander thesis
$wb = $excel.ActiveWorkbook

# --- Sheet "500 bar" ---
$ws = $wb.Worksheets.Item("500 bar")
$ws.Activate()
$ws.Range("B9").Value = 0.04
$ws.Range("B10").Select()

# --- Sheet "LH2" ---
$ws = $wb.Worksheets.Item("LH2")
$ws.Activate()
$ws.Range("B6").Value = 0.08
$ws.Range("B6").Select()

# --- Sheet "LOHC_load" (no data changes) ---
$ws = $wb.Worksheets.Item("LOHC_load")

# --- Sheet "NH3_load" ---
$ws = $wb.Worksheets.Item("NH3_load")
$ws.Activate()
$ws.Range("B2").Value = 2.8090000000000002
$ws.Range("B4").Value = 0.75717000000000001
$ws.Range("B2").Select()

# --- Sheet "NH3_unload" ---
$ws = $wb.Worksheets.Item("NH3_unload")
$ws.Activate()
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 4.2
$ws.Range("B2").Select()

# --- Sheet "LOHC_unload" becomes the active tab, selection on B3 ---
$ws = $wb.Worksheets.Item("LOHC_unload")
$ws.Activate()
$ws.Range("B3").Select()
